$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "2024-10-09 00:00:00"
$ws.Range("B69").Value = 77000
$ws.Range("C69").Value = 10870.64
$ws.Range("D69").Value = 9620.030000000001
$ws.Range("E69").Value = 7.0672
